# Auto-generated edit script applying "gh-pages output" update diff
# Updates "想去人数" (want-to-go count, column F) across all 4 sheets,
# and marks row 2 of sheet "展览" as "不可售" (no longer purchasable) in column G.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 41955
$ws.Range("G2").Value = "不可售"
$ws.Range("F3").Value = 23
$ws.Range("F4").Value = 23
$ws.Range("F5").Value = 9589
$ws.Range("F7").Value = 878
$ws.Range("F8").Value = 906
$ws.Range("F9").Value = 736
$ws.Range("F10").Value = 216
$ws.Range("F12").Value = 301
$ws.Range("F13").Value = 913
$ws.Range("F15").Value = 127
$ws.Range("F16").Value = 739
$ws.Range("F17").Value = 319
$ws.Range("F18").Value = 1421
$ws.Range("F20").Value = 671
$ws.Range("F21").Value = 706
$ws.Range("F22").Value = 466
$ws.Range("F23").Value = 689
$ws.Range("F24").Value = 741
$ws.Range("F28").Value = 504
$ws.Range("F29").Value = 531
$ws.Range("F31").Value = 244
$ws.Range("F32").Value = 933
$ws.Range("F35").Value = 99
$ws.Range("F37").Value = 146
$ws.Range("F38").Value = 398
$ws.Range("F39").Value = 1282
$ws.Range("F41").Value = 1247
$ws.Range("F42").Value = 377
$ws.Range("F43").Value = 96
$ws.Range("F45").Value = 34
$ws.Range("F46").Value = 33
$ws.Range("F47").Value = 46

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 4449
$ws.Range("F7").Value = 334
$ws.Range("F11").Value = 130
$ws.Range("F12").Value = 10
$ws.Range("F13").Value = 58
$ws.Range("F15").Value = 22

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 2035
$ws.Range("F3").Value = 529
$ws.Range("F4").Value = 408

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 2035
$ws.Range("F3").Value = 529
$ws.Range("F4").Value = 23
$ws.Range("F7").Value = 334
$ws.Range("F8").Value = 23
$ws.Range("F9").Value = 9589
$ws.Range("F11").Value = 878
$ws.Range("F12").Value = 878
$ws.Range("F14").Value = 408
$ws.Range("F15").Value = 906
$ws.Range("F16").Value = 130
$ws.Range("F17").Value = 216
$ws.Range("F18").Value = 301
$ws.Range("F19").Value = 913
$ws.Range("F20").Value = 10
$ws.Range("F22").Value = 58
$ws.Range("F23").Value = 739
$ws.Range("F24").Value = 319
$ws.Range("F25").Value = 1421
$ws.Range("F26").Value = 671
$ws.Range("F27").Value = 706
$ws.Range("F28").Value = 466
$ws.Range("F29").Value = 689
$ws.Range("F30").Value = 741
$ws.Range("F33").Value = 504
$ws.Range("F35").Value = 244
$ws.Range("F36").Value = 933
$ws.Range("F40").Value = 99
$ws.Range("F42").Value = 398
$ws.Range("F43").Value = 1247
$ws.Range("F44").Value = 377
$ws.Range("F45").Value = 96
$ws.Range("F47").Value = 34
$ws.Range("F48").Value = 46
